$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 4
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4
$ws.Range("I3").Value = 14
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 6

$ws.Range("B10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4
$ws.Range("I10").Value = 5

$ws.Range("B11").Value = 2
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 9
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 2

$ws.Range("B12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 0

$ws.Range("B14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("M14").Value = 0

$ws.Range("B16").Value = 3
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 3
$ws.Range("K16").Value = 4
$ws.Range("M16").Value = 5

$ws.Range("B17").Value = 11
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 12
$ws.Range("I17").Value = 31
$ws.Range("K17").Value = 8
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = 13
